$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (5th column), shifting the existing
# "Typist"..."Status" columns one place to the right.
$ws.Range("E1:E3").EntireColumn.Insert()

# New column E header + data
$ws.Range("E1").Value = "Client"
$ws.Range("E2").Value = "Accurate"
$ws.Range("E3").Value = "Accurate"

# The inserted column copied the border formatting of column D; the new
# "Client" data cells should be unbordered with just a font color override.
$ws.Range("E2:E3").Borders.LineStyle = -4142
$ws.Range("E2:E3").Font.Color = 0x000000

# New column width (matches column D's width)
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Update selection to match the recorded cursor position in the diff
$ws.Range("H12").Select()
